$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$L1 = $ws.Range("L1")
$Q1 = $ws.Range("Q1")
$R1 = $ws.Range("R1")

$L1.Copy() | Out-Null
$Q1.PasteSpecial(-4122) | Out-Null
$Q1.Value = "drop"

$L1.Copy() | Out-Null
$R1.PasteSpecial(-4122) | Out-Null
$R1.Value = "dropExplanation"
$R1.Borders.Item(7).LineStyle = -4142

# --- Data rows (rows 2-41) ---
for ($r = 2; $r -le 41; $r++) {
    $Lsrc = $ws.Range("L$r")
    $Qcell = $ws.Range("Q$r")
    $Rcell = $ws.Range("R$r")

    $Lsrc.Copy() | Out-Null
    $Qcell.PasteSpecial(-4122) | Out-Null
    $Qcell.Value = $false

    $Lsrc.Copy() | Out-Null
    $Rcell.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# --- Selection / active cell ---
$ws.Range("Q1:R41").Select() | Out-Null
